$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly refresh of fruit/vegetable price data: row values (Fecha, Volumen, Precios) are
# reshuffled across the existing rows 2-21 per the new weekly source extract.

$ws.Range("D2").Value = 44498
$ws.Range("J2").Value = 40
$ws.Range("K2:M2").Value = 4000
$ws.Range("P2").Value = 4000
$ws.Range("D3").Value = 44301
$ws.Range("J3").Value = 40
$ws.Range("K3:M3").Value = 3000
$ws.Range("P3").Value = 3000
$ws.Range("D4").Value = 44312
$ws.Range("J4").Value = 50
$ws.Range("K4:M4").Value = 4000
$ws.Range("P4").Value = 4000
$ws.Range("D5").Value = 44656
$ws.Range("J5").Value = 85
$ws.Range("K5:M5").Value = 5000
$ws.Range("P5").Value = 5000
$ws.Range("D6").Value = 44680
$ws.Range("J6").Value = 20
$ws.Range("K6:M6").Value = 5000
$ws.Range("P6").Value = 5000
$ws.Range("D7").Value = 44365
$ws.Range("J7").Value = 55
$ws.Range("K7:M7").Value = 5000
$ws.Range("P7").Value = 5000
$ws.Range("D8").Value = 44280
$ws.Range("J8").Value = 55
$ws.Range("K8:M8").Value = 4000
$ws.Range("P8").Value = 4000
$ws.Range("D9").Value = 44176
$ws.Range("J9").Value = 10
$ws.Range("K9:M9").Value = 4000
$ws.Range("P9").Value = 4000
$ws.Range("D10").Value = 44390
$ws.Range("J10").Value = 55
$ws.Range("K10:M10").Value = 6000
$ws.Range("P10").Value = 6000
$ws.Range("D11").Value = 44504
$ws.Range("J11").Value = 55
$ws.Range("K11:M11").Value = 4000
$ws.Range("P11").Value = 4000
$ws.Range("D12").Value = 44316
$ws.Range("J12").Value = 20
$ws.Range("K12:M12").Value = 4000
$ws.Range("P12").Value = 4000
$ws.Range("D13").Value = 44509
$ws.Range("J13").Value = 20
$ws.Range("K13:M13").Value = 4000
$ws.Range("P13").Value = 4000
$ws.Range("D14").Value = 44649
$ws.Range("J14").Value = 20
$ws.Range("K14:M14").Value = 5000
$ws.Range("P14").Value = 5000
$ws.Range("D15").Value = 44315
$ws.Range("J15").Value = 40
$ws.Range("K15:M15").Value = 4000
$ws.Range("P15").Value = 4000
$ws.Range("D16").Value = 44259
$ws.Range("J16").Value = 30
$ws.Range("K16:M16").Value = 4000
$ws.Range("P16").Value = 4000
$ws.Range("D17").Value = 44508
$ws.Range("J17").Value = 30
$ws.Range("K17:M17").Value = 4000
$ws.Range("P17").Value = 4000
# Row 18 unchanged
$ws.Range("D19").Value = 44497
$ws.Range("J19").Value = 20
$ws.Range("K19:M19").Value = 4000
$ws.Range("P19").Value = 4000
$ws.Range("D20").Value = 44313
$ws.Range("J20").Value = 20
$ws.Range("K20:M20").Value = 4000
$ws.Range("P20").Value = 4000
$ws.Range("D21").Value = 44291
$ws.Range("J21").Value = 35
$ws.Range("K21:M21").Value = 4000
$ws.Range("P21").Value = 4000
